$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# New "answer_options" (column G) text for the Likert_7 rows (strongly disagree .. strongly agree)
$likert7 = "strongly disagree;  disagree;  somewhat disagree;  neither agree or disagree;  somewhat agree;  agree;  strongly agree"

# New "answer_options" (column G) text for the Likert_6 rows (knowledge scale)
$likert6 = "I need a lot of additional knowledge about the topic;  I need some additional knowledge about the topic ;  I need a little additional knowledge about the topic ;  I have some knowledge about the topic ;  I have good knowledge about the topic ;  I have strong knowledge about the topic"

# Rows that use the Likert_7 answer options (column D already contains "Likert_7";
# row 143 is skipped because it already has a G value in the source file)
$ws.Range("G84:G142").Value = $likert7
$ws.Range("G144:G146").Value = $likert7
$ws.Range("G177:G229").Value = $likert7

# Rows that use the Likert_6 answer options (column D already contains "Likert_6")
$ws.Range("G230:G260").Value = $likert6

# Restore the view state recorded in the saved workbook (window position/size,
# top-left visible cell, and current selection)
$window = $excel.ActiveWindow
$window.ScrollColumn = 3
$window.ScrollRow = 1
[void]$ws.Range("F15").Select()
